# Insert a new data row at row 511 (pushing the existing rows 511-599 down
# to 512-600), then populate the newly inserted row with the new record's
# values. Columns that repeat the pattern of the row above (A,B,C,E,F,G,H,I,
# N,O,Q,R) are filled in to match; D (Fecha), J (Volumen), K/L/M (Precio
# min/max/promedio) and P (Precio $/Kg) carry the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("511:511").Insert()

$ws.Cells.Item(511, 1).Value = 5
$ws.Cells.Item(511, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(511, 3).Value = "Maule"
$ws.Cells.Item(511, 4).Value = 45180
$ws.Cells.Item(511, 5).Value = 7
$ws.Cells.Item(511, 6).Value = 100114013
$ws.Cells.Item(511, 7).Value = "Zanahoria"
$ws.Cells.Item(511, 8).Value = "Sin especificar"
$ws.Cells.Item(511, 9).Value = "Primera"
$ws.Cells.Item(511, 10).Value = 700
$ws.Cells.Item(511, 11).Value = 5000
$ws.Cells.Item(511, 12).Value = 5000
$ws.Cells.Item(511, 13).Value = 5000
$ws.Cells.Item(511, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(511, 15).Value = "Región de Ñuble"
$ws.Cells.Item(511, 16).Value = 250
$ws.Cells.Item(511, 17).Value = 20
$ws.Cells.Item(511, 18).Value = "Hortaliza"
